# "Add condition and more event 1"
# Inserts the canSignUp condition event plus the signUp/signUpStart/signUpDialog/
# signUpWindow/signUpFailed event chain into the eventAction sheet, and updates
# the existing "signUp" row to reference the new eventList/canSignUp condition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing "close" row (old row 13) stays as row 12's content; make room for
# five brand-new rows (new rows 13-17) between it and the "signUp" row (which
# becomes row 18) by inserting five blank rows after row 12.
$ws.Range("A13:A17").EntireRow.Insert()

# Row 12: close (unchanged data, now sitting one row higher)
$ws.Range("A12").Value = "close"
$ws.Range("B12").Value = "关闭窗口"
$ws.Range("C12").Value = "close"
$ws.Range("D12").Value = ";"
$ws.Range("B12").Font.Name = "宋体"
$ws.Range("C12").Font.Name = "宋体"

# Row 13: canSignUp condition event (new). Column C ("condition") was left in
# the default font - unlike the other new rows - matching the source edit.
$ws.Range("A13").Value = "canSignUp"
$ws.Range("B13").Value = "条件分歧"
$ws.Range("C13").Value = "condition"
$ws.Range("D13").Value = "signUpStart;signUpFailed"
$ws.Range("B13").Font.Name = "宋体"
$ws.Range("C13").ClearFormats()

# Row 14: signUpStart (new)
$ws.Range("A14").Value = "signUpStart"
$ws.Range("B14").Value = "签约开始"
$ws.Range("C14").Value = "eventList"
$ws.Range("D14").Value = "signUpDialog;signUpWindow"
$ws.Range("B14").Font.Name = "宋体"
$ws.Range("C14").Font.Name = "宋体"

# Row 15: signUpDialog (new)
$ws.Range("A15").Value = "signUpDialog"
$ws.Range("B15").Value = "签约对话"
$ws.Range("C15").Value = "dialog"
$ws.Range("D15").Value = "dialog_signup_dialog"
$ws.Range("B15").Font.Name = "宋体"
$ws.Range("C15").Font.Name = "宋体"

# Row 16: signUpWindow (new)
$ws.Range("A16").Value = "signUpWindow"
$ws.Range("B16").Value = "签约窗口"
$ws.Range("C16").Value = "window"
$ws.Range("D16").Value = ";"
$ws.Range("B16").Font.Name = "宋体"
$ws.Range("C16").Font.Name = "宋体"

# Row 17: signUpFailed (new)
$ws.Range("A17").Value = "signUpFailed"
$ws.Range("B17").Value = "签约失败对话"
$ws.Range("C17").Value = "dialog"
$ws.Range("D17").Value = "dialog_signup_failure_full"
$ws.Range("B17").Font.Name = "宋体"
$ws.Range("C17").Font.Name = "宋体"

# Row 18: signUp (was row 12) - now routes into the eventList/canSignUp condition
$ws.Range("A18").Value = "signUp"
$ws.Range("B18").Value = "签约"
$ws.Range("C18").Value = "eventList"
$ws.Range("D18").Value = "canSignUp"
$ws.Range("B18").Font.Name = "宋体"
$ws.Range("C18").Font.Name = "宋体"

# Widen column B to fit the new Chinese labels, matching the author's manual
# best-fit resize after typing the new text.
$ws.Columns("B").ColumnWidth = 12.3

# Leave the cursor on the last edited cell, as in the authored workbook.
$ws.Range("D16").Select()
